$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Set cell values in the same order the shared-string table was originally
# built, so new entries land at the expected shared-string indices:
#   31 CollisionSystem        (A11)
#   32 Logic for handling collisions (B11)
#   33 CollusionSystem        (A30)
#   34 ObstacleComponent      (H20 / G39)
#   35 HudRenderSystem        (A12)
#   36 Logic for rendering HUD (B12)
#   37 HudRenderingSystem     (A31)
# ---------------------------------------------------------------------------

# 1. New data: CollisionSystem description row (11)
$ws.Range("A11").Value = "CollisionSystem"
$ws.Range("B11").Value = "Logic for handling collisions"

# 2. New system row in the COMPONENT/SYSTEM matrix (row 30)
$ws.Range("A30").Value = "CollusionSystem"

# 3. New column header "ObstacleComponent" in the two matrices (H20 and G39)
$ws.Range("H20").Value = "ObstacleComponent"
$ws.Range("G39").Value = "ObstacleComponent"

# 4. New data: HudRenderSystem description row (12)
$ws.Range("A12").Value = "HudRenderSystem"
$ws.Range("B12").Value = "Logic for rendering HUD"

# 5. New system row in the COMPONENT/SYSTEM matrix (row 31)
$ws.Range("A31").Value = "HudRenderingSystem"

# ---------------------------------------------------------------------------
# Widen column G to fit the new header
# ---------------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 18.6

# ---------------------------------------------------------------------------
# Mark the component cells the new CollusionSystem row actually touches
# (Bounds/PlayerComponent/ObstacleComponent columns -> B, D, H) as "Good"
# the same way the existing rows do, reusing the existing "Good" style.
# ---------------------------------------------------------------------------
$ws.Range("B23").Copy() | Out-Null
$ws.Range("B30").PasteSpecial(-4122) | Out-Null

$ws.Range("B23").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4122) | Out-Null

$ws.Range("B23").Copy() | Out-Null
$ws.Range("H30").PasteSpecial(-4122) | Out-Null

$ws.Rows.Item(30).RowHeight = 18

# ---------------------------------------------------------------------------
# Mark the ObstacleComponent column for the ObstacleSpawnSystem row (41)
# as "Good" as well
# ---------------------------------------------------------------------------
$ws.Range("B23").Copy() | Out-Null
$ws.Range("G41").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 6. Update the view state: scroll position + active selection
# ---------------------------------------------------------------------------
$ws.Range("B31").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1

Write-Host "edit applied"
